$wb = $excel.ActiveWorkbook

# --- Update the conversion text on "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.1 = 15748.2 pesos`n✅ 15748.2 pesos = 4.07 = 935.7 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 244
$wsTasas.Range("O10").Value = 3842.56
$wsTasas.Range("N12").Value = 3871
$wsTasas.Range("O12").Value = 230
